# Update TPM-derived NATMI ligand-receptor metrics in Sheet1 (F2-F2rl3 pair)
# Columns E..T (5..20) hold recomputed values for rows 2..10 after the
# underlying TPM data was refreshed ("update scripts wuth new tpm").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 5).Value = 1
$ws.Cells.Item(2, 6).Value = 0.3333333333333333
$ws.Cells.Item(2, 7).Value = 0.04339299999999999
$ws.Cells.Item(2, 8).Value = 0.130179
$ws.Cells.Item(2, 9).Value = 0.0698021577815419
$ws.Cells.Item(2, 10).Value = 0.0698021577815419
$ws.Cells.Item(2, 13).Value = 1.600389
$ws.Cells.Item(2, 14).Value = 4.801167
$ws.Cells.Item(2, 15).Value = 0.6625795902775726
$ws.Cells.Item(2, 16).Value = 0.6625795902775727
$ws.Cells.Item(2, 17).Value = 0.06944567987699998
$ws.Cells.Item(2, 18).Value = 0.6250111188929999
$ws.Cells.Item(2, 19).Value = 0.04624948510338452
$ws.Cells.Item(2, 20).Value = 0.04624948510338452
$ws.Cells.Item(3, 5).Value = 1
$ws.Cells.Item(3, 6).Value = 0.3333333333333333
$ws.Cells.Item(3, 7).Value = 0.04339299999999999
$ws.Cells.Item(3, 8).Value = 0.130179
$ws.Cells.Item(3, 9).Value = 0.0698021577815419
$ws.Cells.Item(3, 10).Value = 0.0698021577815419
$ws.Cells.Item(3, 15).Value = 0.253652617229451
$ws.Cells.Item(3, 16).Value = 0.253652617229451
$ws.Cells.Item(3, 17).Value = 0.02658560377433333
$ws.Cells.Item(3, 18).Value = 0.239270433969
$ws.Cells.Item(3, 19).Value = 0.01770550000955119
$ws.Cells.Item(3, 20).Value = 0.0177055000095512
$ws.Cells.Item(4, 5).Value = 1
$ws.Cells.Item(4, 6).Value = 0.3333333333333333
$ws.Cells.Item(4, 7).Value = 0.04339299999999999
$ws.Cells.Item(4, 8).Value = 0.130179
$ws.Cells.Item(4, 9).Value = 0.0698021577815419
$ws.Cells.Item(4, 10).Value = 0.0698021577815419
$ws.Cells.Item(4, 14).Value = 0.606996
$ws.Cells.Item(4, 15).Value = 0.08376779249297628
$ws.Cells.Item(4, 16).Value = 0.0837677924929763
$ws.Cells.Item(4, 17).Value = 0.008779792475999998
$ws.Cells.Item(4, 18).Value = 0.079018132284
$ws.Cells.Item(4, 19).Value = 0.005847172668606192
$ws.Cells.Item(4, 20).Value = 0.005847172668606193
$ws.Cells.Item(5, 9).Value = 0.6764796878879081
$ws.Cells.Item(5, 10).Value = 0.6764796878879081
$ws.Cells.Item(5, 13).Value = 1.600389
$ws.Cells.Item(5, 14).Value = 4.801167
$ws.Cells.Item(5, 15).Value = 0.6625795902775726
$ws.Cells.Item(5, 16).Value = 0.6625795902775727
$ws.Cells.Item(5, 17).Value = 0.6730249227449999
$ws.Cells.Item(5, 18).Value = 6.057224304704999
$ws.Cells.Item(5, 19).Value = 0.4482216344318704
$ws.Cells.Item(5, 20).Value = 0.4482216344318705
$ws.Cells.Item(6, 9).Value = 0.6764796878879081
$ws.Cells.Item(6, 10).Value = 0.6764796878879081
$ws.Cells.Item(6, 15).Value = 0.253652617229451
$ws.Cells.Item(6, 16).Value = 0.253652617229451
$ws.Cells.Item(6, 19).Value = 0.17159084333533
$ws.Cells.Item(6, 20).Value = 0.1715908433353301
$ws.Cells.Item(7, 9).Value = 0.6764796878879081
$ws.Cells.Item(7, 10).Value = 0.6764796878879081
$ws.Cells.Item(7, 14).Value = 0.606996
$ws.Cells.Item(7, 15).Value = 0.08376779249297628
$ws.Cells.Item(7, 16).Value = 0.0837677924929763
$ws.Cells.Item(7, 17).Value = 0.08508836205999999
$ws.Cells.Item(7, 18).Value = 0.76579525854
$ws.Cells.Item(7, 19).Value = 0.05666721012070765
$ws.Cells.Item(7, 20).Value = 0.05666721012070767
$ws.Cells.Item(8, 9).Value = 0.2537181543305499
$ws.Cells.Item(8, 10).Value = 0.2537181543305499
$ws.Cells.Item(8, 13).Value = 1.600389
$ws.Cells.Item(8, 14).Value = 4.801167
$ws.Cells.Item(8, 15).Value = 0.6625795902775726
$ws.Cells.Item(8, 16).Value = 0.6625795902775727
$ws.Cells.Item(8, 17).Value = 0.252422421951
$ws.Cells.Item(8, 18).Value = 2.271801797559
$ws.Cells.Item(8, 19).Value = 0.1681084707423177
$ws.Cells.Item(8, 20).Value = 0.1681084707423177
$ws.Cells.Item(9, 9).Value = 0.2537181543305499
$ws.Cells.Item(9, 10).Value = 0.2537181543305499
$ws.Cells.Item(9, 15).Value = 0.253652617229451
$ws.Cells.Item(9, 16).Value = 0.253652617229451
$ws.Cells.Item(9, 18).Value = 0.8697045309470001
$ws.Cells.Item(9, 19).Value = 0.06435627388456976
$ws.Cells.Item(9, 20).Value = 0.06435627388456977
$ws.Cells.Item(10, 9).Value = 0.2537181543305499
$ws.Cells.Item(10, 10).Value = 0.2537181543305499
$ws.Cells.Item(10, 14).Value = 0.606996
$ws.Cells.Item(10, 15).Value = 0.08376779249297628
$ws.Cells.Item(10, 16).Value = 0.0837677924929763
$ws.Cells.Item(10, 17).Value = 0.031912949588
$ws.Cells.Item(10, 18).Value = 0.287216546292
$ws.Cells.Item(10, 19).Value = 0.02125340970366244
$ws.Cells.Item(10, 20).Value = 0.02125340970366245
